$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.666.68'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '1.593.48'
$ws.Range('E3').Value = '  -1.60%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''211.73'
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').Value = '''0.0618'
$ws.Range('E8').Value = '  -1.47%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = '''0.246'
$ws.Range('E9').Value = '  -2.60%  '
$ws.Range('D10').Value = '''19.63'
$ws.Range('E10').Value = '  -2.18%  '
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('D12').Value = '1.816.30'
$ws.Range('E12').Value = '  -1.60%  '
$ws.Range('D13').Value = '1.605.26'
$ws.Range('E13').Value = '  -0.92%  '
$ws.Range('E14').Value = '  -2.69%  '
$ws.Range('E15').Value = '  -3.16%  '
$ws.Range('D16').Value = '''65.34'
$ws.Range('E16').Value = '  +1.02%  '
$ws.Range('D17').Value = '26.628.40'
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('E18').Value = '  -2.73%  '
$ws.Range('D19').Value = '''209.28'
$ws.Range('E19').Value = '  -2.47%  '
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').Value = '''6.71'
$ws.Range('E21').Value = '  -2.30%  '
$ws.Range('E22').Value = '  -2.36%  '
$ws.Range('D23').Value = '''2.32'
$ws.Range('E23').Value = '  -3.46%  '
$ws.Range('D24').Value = '''8.90'
$ws.Range('E24').Value = '  -1.95%  '
$ws.Range('D25').Value = '''145.81'
$ws.Range('E25').Value = '  -1.67%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').Value = '''7.18'
$ws.Range('E27').Value = '  -3.38%  '
$ws.Range('E28').Value = '  -1.91%  '
$ws.Range('D29').Value = '''15.34'
$ws.Range('E29').Value = '  -1.42%  '
$ws.Range('E30').Value = '  -2.41%  '
$ws.Range('E31').Value = '  -1.65%  '
$ws.Range('D32').Value = '''3.23'
$ws.Range('E32').Value = '  -4.17%  '
$ws.Range('D33').Value = '''0.665'
$ws.Range('E33').Value = '  -12.38%  '
$ws.Range('E34').Value = '  -3.28%  '
$ws.Range('D35').Value = '1.306.86'
$ws.Range('E35').Value = '  -2.79%  '
$ws.Range('D36').Value = '''2.45'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('E37').Value = '  -4.89%  '
$ws.Range('D38').Value = '''0.0172'
$ws.Range('E38').Value = '  -4.05%  '
$ws.Range('E39').Value = '  -2.83%  '
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('E41').Value = '  -1.18%  '
$ws.Range('D42').Value = '''5.36'
$ws.Range('E42').Value = '  +1.04%  '
$ws.Range('E43').Value = '  -1.87%  '
$ws.Range('D44').Value = '''63.18'
$ws.Range('E44').Value = '  -3.28%  '
$ws.Range('D45').Value = '1.729.25'
$ws.Range('E45').Value = '  -1.39%  '
$ws.Range('D46').Value = '''89.08'
$ws.Range('E46').Value = '  -1.06%  '
$ws.Range('D47').Value = '''1.60'
$ws.Range('E47').Value = '  -1.94%  '
$ws.Range('D48').Value = '''0.819'
$ws.Range('E48').Value = '  -8.77%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '''0.0978'
$ws.Range('E49').Value = '  -2.22%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.0503'
$ws.Range('E50').Value = '  -2.57%  '
$ws.Range('D51').Value = '''7.51'
$ws.Range('E51').Value = '  -1.90%  '
